# Apply the change described in the commit: "added data for username and password."
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "Script1" to "ValidLogin"
$ws.Name = "ValidLogin"

# Add new header + value for Password column (B)
$ws.Range("B1").Value = "Password"
$ws.Range("B2").Value = "manager"

# Move selection to B3 (matches post-edit selection in the diff)
$ws.Range("B3").Select()
